$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the TISG_PDR_G (B), fcs (C) and need_to_buy_MW (F) forecast columns
# for rows 2-7. Re-applying the default "Normal" style after clearing keeps
# the now-blank cells present in the sheet (rather than dropping them
# entirely), matching the target workbook which still carries empty <c/>
# placeholders for these cells.
$ws.Range("B2:C7").ClearContents()
$ws.Range("B2:C7").Style = "Normal"

$ws.Range("F2:F7").ClearContents()
$ws.Range("F2:F7").Style = "Normal"

# Update buy_BEE_MWH (D2) from 2376 to 5016
$ws.Range("D2").Value = 5016
